$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1750183612496/A5MTrC4t3W-camera-video-segment-driver-1750183614996.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSICANEIFG%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T030004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK7%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCwKeVwP%2FGu7v4ZCQ8vbQkTsA19DXNXA1AfRsB1cyFf8gIgRGPjutNZpYth1J%2B0r9SbqcsL%2BKhYB2sCX%2FJx8WbQBwEq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDIDKaAq0mTZeX8o9liq6AwSOtY6k5FVe%2B7loK%2FPc9au3nnNY2lFAPUvCnIeWr35TnIYRmiOzZMM8xZo2nF0uDzzgLrcLOm7unuR4Y6od8gHJvnnFq4wDhpvyaNu3pMWdTta88x%2FDlhKIWQCnSsyZzfSPigIWWuPxrfwvj3NC5m9qLtXSscaKGpBZ765aqKZMlDpppcfskp5zinbDXOkTvTSqa6XEtJ7d8tljPaV8ERQTIBdn02CJtVlvX0OYZz6CAyY2Dyv6pkjiHN6eEKxm0N5%2FARy6k0vPwg5IjzgQL7yvSvpClrSHi9CtjBKsMHhfTC5%2BNgZw14a8NYjyLAJPD7vDXzYc%2BPRDXK2ZAWnCNlliD4HykRq5FHgamNVF9etjNWTLGcl7PTvZpV9z5U%2BgZ3tH2YJtsbPVDfhL7KkuIhg3GWcIimtesSSCgbFsF7Z3WcKjKnOuTvgR7Jnve%2Bzau9jIdb2mHy5TXg%2FDeFn8SnWUPd1SYgm78pfepln8ymy04qCKBrqyj6YfPdix8AA5qTeQZgB1VGmV6TnmfIHlDQx85JHOMPX3uAIe%2FHeBCP2U34iq%2BvHGfJJ1aAyCimvrC0ACUHB1uxZuq1kw3NrMwgY6pQE0Zu9E1t%2BEvap9tP34BFqVD3%2Bf78Gcc1vd%2FPmUTOw6Dk4r1wiRZiKpPUPKHuYu5IyAaEokgVeBHR08lVWFez2pqzjf2gaVKMwEenIZgDUsWpBsHMKqGBkoDdfcuqQ6n0eWDyCB1w%2BUjNEZRGpNIG4QHQdaiXVmCfsz5a2Tg0gHM0C3xtSfYL3XgKomrtsKYcCoknLb7E6R0Gk%2ByDP4aunKka%2BocdI%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2011%3A00%3A04%20GMT&X-Amz-Signature=b2f39505c33225577f531e17b6f5c0950f927669741fd2d2e98b13b501bb8b71'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1750183612496/TAdN8NokD3-camera-video-segment-1750183614996.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSICANEIFG%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T030004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK7%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCwKeVwP%2FGu7v4ZCQ8vbQkTsA19DXNXA1AfRsB1cyFf8gIgRGPjutNZpYth1J%2B0r9SbqcsL%2BKhYB2sCX%2FJx8WbQBwEq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDIDKaAq0mTZeX8o9liq6AwSOtY6k5FVe%2B7loK%2FPc9au3nnNY2lFAPUvCnIeWr35TnIYRmiOzZMM8xZo2nF0uDzzgLrcLOm7unuR4Y6od8gHJvnnFq4wDhpvyaNu3pMWdTta88x%2FDlhKIWQCnSsyZzfSPigIWWuPxrfwvj3NC5m9qLtXSscaKGpBZ765aqKZMlDpppcfskp5zinbDXOkTvTSqa6XEtJ7d8tljPaV8ERQTIBdn02CJtVlvX0OYZz6CAyY2Dyv6pkjiHN6eEKxm0N5%2FARy6k0vPwg5IjzgQL7yvSvpClrSHi9CtjBKsMHhfTC5%2BNgZw14a8NYjyLAJPD7vDXzYc%2BPRDXK2ZAWnCNlliD4HykRq5FHgamNVF9etjNWTLGcl7PTvZpV9z5U%2BgZ3tH2YJtsbPVDfhL7KkuIhg3GWcIimtesSSCgbFsF7Z3WcKjKnOuTvgR7Jnve%2Bzau9jIdb2mHy5TXg%2FDeFn8SnWUPd1SYgm78pfepln8ymy04qCKBrqyj6YfPdix8AA5qTeQZgB1VGmV6TnmfIHlDQx85JHOMPX3uAIe%2FHeBCP2U34iq%2BvHGfJJ1aAyCimvrC0ACUHB1uxZuq1kw3NrMwgY6pQE0Zu9E1t%2BEvap9tP34BFqVD3%2Bf78Gcc1vd%2FPmUTOw6Dk4r1wiRZiKpPUPKHuYu5IyAaEokgVeBHR08lVWFez2pqzjf2gaVKMwEenIZgDUsWpBsHMKqGBkoDdfcuqQ6n0eWDyCB1w%2BUjNEZRGpNIG4QHQdaiXVmCfsz5a2Tg0gHM0C3xtSfYL3XgKomrtsKYcCoknLb7E6R0Gk%2ByDP4aunKka%2BocdI%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2011%3A00%3A04%20GMT&X-Amz-Signature=571b4ebd6be8bea33e1965adcae718e2e549e7d2573b2a2d3fc7d6ec77e61d80'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1750182961683/iBjUYEdqdu-camera-video-segment-driver-1750182964183.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSICANEIFG%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T030004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK7%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCwKeVwP%2FGu7v4ZCQ8vbQkTsA19DXNXA1AfRsB1cyFf8gIgRGPjutNZpYth1J%2B0r9SbqcsL%2BKhYB2sCX%2FJx8WbQBwEq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDIDKaAq0mTZeX8o9liq6AwSOtY6k5FVe%2B7loK%2FPc9au3nnNY2lFAPUvCnIeWr35TnIYRmiOzZMM8xZo2nF0uDzzgLrcLOm7unuR4Y6od8gHJvnnFq4wDhpvyaNu3pMWdTta88x%2FDlhKIWQCnSsyZzfSPigIWWuPxrfwvj3NC5m9qLtXSscaKGpBZ765aqKZMlDpppcfskp5zinbDXOkTvTSqa6XEtJ7d8tljPaV8ERQTIBdn02CJtVlvX0OYZz6CAyY2Dyv6pkjiHN6eEKxm0N5%2FARy6k0vPwg5IjzgQL7yvSvpClrSHi9CtjBKsMHhfTC5%2BNgZw14a8NYjyLAJPD7vDXzYc%2BPRDXK2ZAWnCNlliD4HykRq5FHgamNVF9etjNWTLGcl7PTvZpV9z5U%2BgZ3tH2YJtsbPVDfhL7KkuIhg3GWcIimtesSSCgbFsF7Z3WcKjKnOuTvgR7Jnve%2Bzau9jIdb2mHy5TXg%2FDeFn8SnWUPd1SYgm78pfepln8ymy04qCKBrqyj6YfPdix8AA5qTeQZgB1VGmV6TnmfIHlDQx85JHOMPX3uAIe%2FHeBCP2U34iq%2BvHGfJJ1aAyCimvrC0ACUHB1uxZuq1kw3NrMwgY6pQE0Zu9E1t%2BEvap9tP34BFqVD3%2Bf78Gcc1vd%2FPmUTOw6Dk4r1wiRZiKpPUPKHuYu5IyAaEokgVeBHR08lVWFez2pqzjf2gaVKMwEenIZgDUsWpBsHMKqGBkoDdfcuqQ6n0eWDyCB1w%2BUjNEZRGpNIG4QHQdaiXVmCfsz5a2Tg0gHM0C3xtSfYL3XgKomrtsKYcCoknLb7E6R0Gk%2ByDP4aunKka%2BocdI%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2011%3A00%3A04%20GMT&X-Amz-Signature=55e752d97977af7c819670b98216e347cf258aefce5e03aa23343b45ef5a03fa'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1750182961683/Foph7ZSROm-camera-video-segment-1750182964183.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSICANEIFG%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T030004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEK7%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCwKeVwP%2FGu7v4ZCQ8vbQkTsA19DXNXA1AfRsB1cyFf8gIgRGPjutNZpYth1J%2B0r9SbqcsL%2BKhYB2sCX%2FJx8WbQBwEq5gMIlv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDIDKaAq0mTZeX8o9liq6AwSOtY6k5FVe%2B7loK%2FPc9au3nnNY2lFAPUvCnIeWr35TnIYRmiOzZMM8xZo2nF0uDzzgLrcLOm7unuR4Y6od8gHJvnnFq4wDhpvyaNu3pMWdTta88x%2FDlhKIWQCnSsyZzfSPigIWWuPxrfwvj3NC5m9qLtXSscaKGpBZ765aqKZMlDpppcfskp5zinbDXOkTvTSqa6XEtJ7d8tljPaV8ERQTIBdn02CJtVlvX0OYZz6CAyY2Dyv6pkjiHN6eEKxm0N5%2FARy6k0vPwg5IjzgQL7yvSvpClrSHi9CtjBKsMHhfTC5%2BNgZw14a8NYjyLAJPD7vDXzYc%2BPRDXK2ZAWnCNlliD4HykRq5FHgamNVF9etjNWTLGcl7PTvZpV9z5U%2BgZ3tH2YJtsbPVDfhL7KkuIhg3GWcIimtesSSCgbFsF7Z3WcKjKnOuTvgR7Jnve%2Bzau9jIdb2mHy5TXg%2FDeFn8SnWUPd1SYgm78pfepln8ymy04qCKBrqyj6YfPdix8AA5qTeQZgB1VGmV6TnmfIHlDQx85JHOMPX3uAIe%2FHeBCP2U34iq%2BvHGfJJ1aAyCimvrC0ACUHB1uxZuq1kw3NrMwgY6pQE0Zu9E1t%2BEvap9tP34BFqVD3%2Bf78Gcc1vd%2FPmUTOw6Dk4r1wiRZiKpPUPKHuYu5IyAaEokgVeBHR08lVWFez2pqzjf2gaVKMwEenIZgDUsWpBsHMKqGBkoDdfcuqQ6n0eWDyCB1w%2BUjNEZRGpNIG4QHQdaiXVmCfsz5a2Tg0gHM0C3xtSfYL3XgKomrtsKYcCoknLb7E6R0Gk%2ByDP4aunKka%2BocdI%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2011%3A00%3A04%20GMT&X-Amz-Signature=fd5e12dc54895c71afc7133740b62a26f7cdceb44fa93f63648f3de6f1c38cfd'
